# Rewrites Sheet1 rows 2-13 with refreshed TPM-based NATMI ligand/receptor
# statistics, and adds three new rows (11-13) for the 'Resolving-Mac' sending
# cluster (a population not present in the previous run).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rowData = @{
    2 = @('ECs', 'Colq', 'Musk', 'ECs', 1, 0.3333333333333333, 0.004118333333333333, 0.012355, 0.005796697547980517, 0.005796697547980516, 2, 0.6666666666666666, 0.4200680000000001, 1.260204, 0.01461602726853518, 0.01461602726853518, 0.001729980046666667, 0.01556982042, 0.00008472468942873425, 0.00008472468942873423)
    3 = @('ECs', 'Colq', 'Musk', 'FAPs', 1, 0.3333333333333333, 0.004118333333333333, 0.012355, 0.005796697547980517, 0.005796697547980516, 3, 1, 9.389971666666668, 28.169915, 0.3267187263271011, 0.3267187263271011, 0.03867103331388889, 0.348039299825, 0.001893889639779625, 0.001893889639779624)
    4 = @('ECs', 'Colq', 'Musk', 'MuSCs', 1, 0.3333333333333333, 0.004118333333333333, 0.012355, 0.005796697547980517, 0.005796697547980516, 3, 1, 18.93019133333333, 56.79057399999999, 0.6586652464043636, 0.6586652464043636, 0.07796083797444443, 0.7016475417699999, 0.003818083218772157, 0.003818083218772157)
    5 = @('FAPs', 'Colq', 'Musk', 'ECs', 3, 1, 0.09240733333333333, 0.277222, 0.1300665388625054, 0.1300665388625054, 2, 0.6666666666666666, 0.4200680000000001, 1.260204, 0.01461602726853518, 0.01461602726853518, 0.03881736369866667, 0.349356273288, 0.00190105607873837, 0.00190105607873837)
    6 = @('FAPs', 'Colq', 'Musk', 'FAPs', 3, 1, 0.09240733333333333, 0.277222, 0.1300665388625054, 0.1300665388625054, 3, 1, 9.389971666666668, 28.169915, 0.3267187263271011, 0.3267187263271011, 0.8677022417922222, 7.80932017613, 0.04249517391493218, 0.04249517391493217)
    7 = @('FAPs', 'Colq', 'Musk', 'MuSCs', 3, 1, 0.09240733333333333, 0.277222, 0.1300665388625054, 0.1300665388625054, 3, 1, 18.93019133333333, 56.79057399999999, 0.6586652464043636, 0.6586652464043636, 1.749288500603111, 15.743596505428, 0.08567030886883488, 0.08567030886883487)
    8 = @('MuSCs', 'Colq', 'Musk', 'ECs', 3, 1, 0.588616, 1.765848, 0.8284975128859814, 0.8284975128859813, 2, 0.6666666666666666, 0.4200680000000001, 1.260204, 0.01461602726853518, 0.01461602726853518, 0.247258745888, 2.225328712992, 0.01210934224025508, 0.01210934224025508)
    9 = @('MuSCs', 'Colq', 'Musk', 'FAPs', 3, 1, 0.588616, 1.765848, 0.8284975128859814, 0.8284975128859813, 3, 1, 9.389971666666668, 28.169915, 0.3267187263271011, 0.3267187263271011, 5.527087562546668, 49.74378806292001, 0.2706856521752789, 0.2706856521752788)
    10 = @('MuSCs', 'Colq', 'Musk', 'MuSCs', 3, 1, 0.588616, 1.765848, 0.8284975128859814, 0.8284975128859813, 3, 1, 18.93019133333333, 56.79057399999999, 0.6586652464043636, 0.6586652464043636, 11.14261350186133, 100.283521516752, 0.5457025184704474, 0.5457025184704473)
    11 = @('Resolving-Mac', 'Colq', 'Musk', 'ECs', 2, 0.6666666666666666, 0.02532033333333333, 0.07596099999999999, 0.03563925070353282, 0.03563925070353281, 2, 0.6666666666666666, 0.4200680000000001, 1.260204, 0.01461602726853518, 0.01461602726853518, 0.01063626178266667, 0.095726356044, 0.0005209042601129972, 0.0005209042601129971)
    12 = @('Resolving-Mac', 'Colq', 'Musk', 'FAPs', 2, 0.6666666666666666, 0.02532033333333333, 0.07596099999999999, 0.03563925070353282, 0.03563925070353281, 3, 1, 9.389971666666668, 28.169915, 0.3267187263271011, 0.3267187263271011, 0.2377572125905555, 2.139814913315, 0.01164401059711049, 0.01164401059711048)
    13 = @('Resolving-Mac', 'Colq', 'Musk', 'MuSCs', 2, 0.6666666666666666, 0.02532033333333333, 0.07596099999999999, 0.03563925070353282, 0.03563925070353281, 3, 1, 18.93019133333333, 56.79057399999999, 0.6586652464043636, 0.6586652464043636, 0.4793187546237776, 4.313868791613999, 0.02347433584630933, 0.02347433584630933)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value2 = $vals[$i]
    }
}

